# Applies the "network log display logic" update to the language sheet:
#  - adds a "search in files" title string right after the existing
#    "search_title" row
#  - adds a "Network Log" inspector title row right before the existing
#    "files" row
#  - appends four new network-connection-state strings at the end of the
#    sheet
#  - refreshes the sheet view (scroll position / selection) and used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "search_title_files" row right after row 17 (search_title) ---
$ws.Rows("18:18").Insert()
$ws.Range("A18").Value = "search_title_files"
$ws.Range("B18").Value = "Search In Files: {0}"

# --- Insert "dataInvestigate_network_log_inspector" row right before the
#     "files" row (now shifted down to row 27) ---
$ws.Rows("27:27").Insert()
$ws.Range("A27").Value = "dataInvestigate_network_log_inspector"
$ws.Range("B27").Value = "Network Log"

# --- Append the new network connection state strings at the bottom ---
$ws.Range("A36").Value = "network_log_state_listening"
$ws.Range("B36").Value = "LISTENING"

$ws.Range("A37").Value = "network_log_state_established"
$ws.Range("B37").Value = "ESTABLISHED"

$ws.Range("A38").Value = "network_log_state_close_wait"
$ws.Range("B38").Value = "CLOSE_WAIT"

$ws.Range("A39").Value = "network_log_state_time_wait"
$ws.Range("B39").Value = "TIME_WAIT"

# --- Update sheet view: scroll position + current selection ---
$null = $ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 6
